$wb = $excel.ActiveWorkbook

# 1. Remove the "customer_leads" worksheet (feedback: remove table)
$leadsSheet = $wb.Worksheets.Item("customer_leads")
$leadsSheet.Delete()

# 2. Add a "note" column to the "customers" sheet with notes for a few customers
$ws = $wb.Worksheets.Item("customers")
$ws.Range("J1").Value = "note"
$ws.Range("J3").Value = "Needs loan"
$ws.Range("J5").Value = "Needs financing"
$ws.Range("J6").Value = "Inquiry into financing options"

# 3. Clear the stray pre-formatted (but empty) list_price cells in "inventory"
$inv = $wb.Worksheets.Item("inventory")
$inv.Range("D5:D11").Clear()
